$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ---
# Values are written in this specific order so the shared-strings table
# gets appended in the same sequence as the source edit (C5, C6, A5,
# D5, D6, E5, E6).
$ws.Cells.Item(5, 3).Value = " How was your expedition?"
$ws.Cells.Item(6, 3).Value = " Did you happen to find\nany treasures?"

$ws.Cells.Item(5, 1).Value = "SCRIPT/G01P03A/um1102.ssb"

$ws.Cells.Item(5, 4).Value = " Как прошла экспедиция?"
$ws.Cells.Item(6, 4).Value = " Вы нашли какие-нибудь\nсокровища?"

$ws.Cells.Item(5, 5).Value = " Ëàë ðñïšìà üëòðåäéøéÿ?"
$ws.Cells.Item(6, 5).Value = " Âú îàšìé ëàëéå-îéáôäû\nòïëñïâéþà?"

$ws.Cells.Item(5, 2).Value = 158
$ws.Cells.Item(6, 2).Value = 161

# --- Formatting: columns A/B use the "wide" style (font size 11, wrap),
# columns C/D/E use the "narrow" style (font size 8, wrap) - matching the
# style already used by rows 2 and 3. Row 6 has no A6 cell (mirrors row 4,
# which also has no leading-column value), so it is skipped.
$wideCells5 = 1, 2
$wideCells6 = 2
$narrowCols = 3, 4, 5

foreach ($c in $wideCells5) {
    $ws.Cells.Item(5, $c).WrapText = $true
    $ws.Cells.Item(5, $c).Font.Size = 11
}
foreach ($c in $wideCells6) {
    $ws.Cells.Item(6, $c).WrapText = $true
    $ws.Cells.Item(6, $c).Font.Size = 11
}

foreach ($c in $narrowCols) {
    $ws.Cells.Item(5, $c).WrapText = $true
    $ws.Cells.Item(5, $c).Font.Size = 8
    $ws.Cells.Item(6, $c).WrapText = $true
    $ws.Cells.Item(6, $c).Font.Size = 8
}

# --- Row heights ---
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 27.6

# --- Selection ---
$ws.Range("D5").Select() | Out-Null
